$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Prepare formatting for the two new rows (16 and 17) by copying the
# format of row 15 (the last existing data row) into them first, so the
# A-column style ("s=1": bold/border/centered) carries over.
$ws.Range("A15").Copy() | Out-Null
$ws.Range("A16:A17").PasteSpecial(-4122) | Out-Null

# Data rows 8-17: (index, name, from_bus, to_bus, in_service)
$data = @(
    @(8,  "line7", 14, 11, $true),
    @(9,  "line8", 16, 9,  $true),
    @(10, "extr1", 5,  12, $true),
    @(11, "extr2", 5,  9,  $true),
    @(12, "extr3", 10, 11, $true),
    @(13, "extr4", 7,  8,  $true),
    @(14, "extr5", 9,  11, $false),
    @(15, "extr6", 7,  11, $true),
    @(16, "extr7", 5,  7,  $true),
    @(17, "extr8", 8,  5,  $true)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $r - 2
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}
